$p = $ppt.ActivePresentation
$sm = $p.SlideMaster
$cs = $sm.ColorScheme
$c1 = $cs.Item(1)
$c1.RGB = 255
